$d = $word.ActiveDocument

function Find-ParagraphByText($doc, $needle) {
    $paras = $doc.Paragraphs
    for ($i = 1; $i -le $paras.Count; $i++) {
        $p = $paras.Item($i)
        # Range.Text includes the trailing paragraph mark (CR, chr 13) -
        # strip it before comparing against a plain-text needle.
        $t = $p.Range.Text.TrimEnd([char]13)
        if ($t -eq $needle) {
            return $p
        }
    }
    return $null
}

# ---------------------------------------------------------------------------
# 1) Swap the two OLE_LINK bookmarks on "Zbus Architecture" (title paragraph):
#    OLE_LINK2/OLE_LINK1 (ids 0/1) -> OLE_LINK1/OLE_LINK2 (ids 0/1)
# ---------------------------------------------------------------------------
$bmA = $d.Bookmarks.Item("OLE_LINK1")
$startA = $bmA.Start
$endA = $bmA.End
$bmB = $d.Bookmarks.Item("OLE_LINK2")
$startB = $bmB.Start
$endB = $bmB.End

$bmA.Delete()
$bmB.Delete()

$d.Bookmarks.Add("OLE_LINK1", $d.Range($startA, $endA))
$d.Bookmarks.Add("OLE_LINK2", $d.Range($startB, $endB))

# ---------------------------------------------------------------------------
# 2) Add a new bookmark OLE_LINK3 around the "zbus URL pattern" section -
#    starting right at the "zbus URL pattern" heading and ending right after
#    the "/rpc/topic/method/param1/param2/.../[?module=xxx]" paragraph.
# ---------------------------------------------------------------------------
$pStart = Find-ParagraphByText $d "zbus URL pattern"
$pEnd = Find-ParagraphByText $d "/rpc/topic/method/param1/param2/…/[?module=xxx]"

$d.Bookmarks.Add("OLE_LINK3", $d.Range($pStart.Range.Start, $pEnd.Range.End))

# ---------------------------------------------------------------------------
# 3) Remove the old _GoBack bookmark (it used to sit right after "/track_pub")
# ---------------------------------------------------------------------------
$oldGoBack = $d.Bookmarks.Item("_GoBack")
$oldGoBack.Delete()

# ---------------------------------------------------------------------------
# 4) Re-create _GoBack around the "Zbus Client Platforms" heading paragraph.
# ---------------------------------------------------------------------------
$pClient = Find-ParagraphByText $d "Zbus Client Platforms"
$d.Bookmarks.Add("_GoBack", $d.Range($pClient.Range.Start, $pClient.Range.End))

Write-Output "done"
